$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "unknown"
$ws.Range("B6").Value = "unknown"
$ws.Range("C6").Value = "XD"

# Force "1" to be stored as text (not a number) without leaving any
# leftover cell-style behind: compute it as text via TEXT(), then copy
# only the resulting value into D6.
$ws.Range("Z1").Formula = '=TEXT(1,"0")'
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("E6").Value = "brak"
$ws.Range("F6").Value = "mateusz.qzera@gmail.com"
